# Auto-applies the enron_gemma3_1b_prompt6 diff: adds three new model-output
# columns (E, F, G) alongside the existing D column, relabels the D1 header,
# and refreshes D2:D4 with the new Private_gemma3:1b outputs while clearing
# D5:D12 (those models no longer have a recorded prediction for those rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header row - relabel D1, add E1/F1/G1
$ws.Range("D1").Value = "Private_gemma3:1b"
$ws.Range("E1").Value = "Private_gemma:2b"
$ws.Range("F1").Value = "Private_llama3.2:3b"
$ws.Range("G1").Value = "Private_mistral"
$ws.Range("E1").Style = $ws.Range("D1").Style
$ws.Range("F1").Style = $ws.Range("D1").Style
$ws.Range("G1").Style = $ws.Range("D1").Style

# Row 2
$v_D2 = @'
“forecast”: “Here is our forecast”
'@
$ws.Range("D2").Value = $v_D2
$v_E2 = @'
Sure, here's the modified sentence with the names, places, and numbers converted to the "key": "value" format:
"Please write a greeting card for "person_name": "Nancy" when she is "person_age": "18" years old and lives in "person_city": "Boston"."
'@
$ws.Range("E2").Value = $v_E2
$v_F2 = @'
Here is the modified sentence:
'Please write a greeting card for "name": "None" when she is "age": "None" years old and lives in "city": "None"'
'@
$ws.Range("F2").Value = $v_F2
$v_G2 = @'
Here is our forecast with placeholders replaced as follows:
   "Location": "New York",
   "ForecastDay": "Thursday",
   "TemperatureMin": "20",
   "TemperatureMax": "30",
   "WeatherCondition": "Sunny"
   Output: 'Here is our forecast for "Location": "New York" on "ForecastDay": "Thursday". The temperature will range from "TemperatureMin": "20" to "TemperatureMax": "30" degrees, and the weather condition will be "WeatherCondition": "Sunny".'
'@
$ws.Range("G2").Value = $v_G2

# Row 3
$v_D3 = @'
“key": "Austin", "activity": "Play golf", "location": "Austin", "activity": "Rent a ski boat and jet ski's", "location": "Austin"
'@
$ws.Range("D3").Value = $v_D3
$v_E3 = @'
Sure, here's the output as requested:
Traveling to have a business meeting takes the fun out of the trip. Especially if you have to prepare a presentation. I would suggest holding the business meetings here then take a trip without any formal business meetings. I would even try and get some honest opinions on whether a trip is even desired or necessary.
'@
$ws.Range("E3").Value = $v_E3
$v_F3 = @'
'Please write a greeting card for "person name": None when she is "age": None years old and lives in "city": None.'
However, since the original text has multiple instances of specific information, we can rephrase the input to better accommodate this. Here's an example with multiple mentions:
 'Please write a greeting card for "name": "John" when he is "age": "18" years old and lives in "city": "None"'
But we also need to address more instances such as business meetings, time etc.
 Let me try rephrasing it. 
'Traveling to have a business meeting takes the fun out of the trip. Especially if you have to prepare a presentation at 2:00 PM on Friday. I would suggest holding the business plan meetings here then take a trip without any formal business meetings. I would even try and get some honest opinions on whether a trip is even desired or necessary on Saturday.
 As far as the business meetings, I think it would be more productive to try and stimulate discussions across the different groups about what is working and what is not. Too often the presenter speaks and the others are quiet just waiting for their turn in a 60 minute meeting format.
   My suggestion for where to go is "location": "Austin". Play golf and rent a ski boat and jet ski's. Flying somewhere takes too much time.'
Output:
 'Please write a greeting card for "person name": None when she is "age": None years old and lives in "city": None.
 Please prepare a presentation at "time": "2:00 PM" on Friday
 Traveling to have a business meeting takes the fun out of the trip. Especially if you have to prepare a presentation at "time": "2:00 PM" on Friday. I would suggest holding the business plan meetings here then take a trip without any formal business meetings. I would even try and get some honest opinions on whether a trip is even desired or necessary on "day": "Saturday".
 As far as the business meetings, I think it would be more productive to try and stimulate discussions across the different groups about what is working and what is not. Too often the presenter speaks and the others are quiet just waiting for their turn in a "meeting duration": "60 minute" meeting format.
   My suggestion for where to go is "location": "Austin". Play golf and rent a ski boat and jet ski's. Flying somewhere takes too much time.'
'@
$ws.Range("F3").Value = $v_F3
$v_G3 = @'
Traveling to have a business meeting takes the fun out of the trip, especially if one has to prepare a "presentation": "None". I would suggest holding the "business plan meetings" here then take a trip without any formal "business meetings". I would even try and get some "opinions" on whether a trip is even desired or necessary.
   As far as the "business meetings", I think it would be more productive to try and stimulate discussions across the different groups about what is working and what is not. Too often the "presenter" speaks and the others are quiet just waiting for their turn.  The meetings might be better if held in a "format": "round table discussion".
   My suggestion for where to go is "location": "Austin". Play golf and rent a "vehicle_1": "ski boat" and "vehicle_2": "jet ski's". Flying somewhere takes too much time.
   In this output, I replaced specific names with "None", changed the presentation and opinions to a generic term "presentation" and "opinions" respectively, and converted the location, golf equipment, and jet-ski to their respective keys.
'@
$ws.Range("G3").Value = $v_G3

# Row 4
$v_D4 = @'
'key': "test"
'way': "successful"
'to': "go"
'here': "None"
'@
$ws.Range("D4").Value = $v_D4
$v_E4 = @'
Sure, here's the output based on the input you provided:
"Please write a greeting card for "name": "test successful. way to go!!!" when she is "age": "18" years old and lives in "city": "Boston".'
'@
$ws.Range("E4").Value = $v_E4
$v_F4 = @'
'Please write a greeting card for "person Names": None when she is "age": "None" years old and lives in "city": "None", but test was successful and "important time": "None".'
However, since all the specific information in your input was 'None', I've kept that value. If there were any other types of information in the sentence (like phone numbers or organization names), they would have been replaced with their respective key-value pairs.
Here's a revised version of the output to better reflect the result:
'Please write a greeting card for "person Names": None when she is "age": None years old and lives in "city": None, but "test": "successful" and there was indeed a "way to go"!!! and also "time": "None".'
Let me know if you need further clarification or have any other questions!
'@
$ws.Range("F4").Value = $v_F4
$v_G4 = @'
None (No sensitive or specific information found in the sentence)
'@
$ws.Range("G4").Value = $v_G4

# Rows 5-12: these model columns have no recorded prediction any more -
# clear D (it previously held a value) and leave E/F/G blank (they never
# had values, so there is nothing to clear there).
$rowsToClearD = 5..12
foreach ($r in $rowsToClearD) {
    $ws.Range("D" + $r).ClearContents()
}

